$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the whole "#include <string>" paragraph and the whole
#    "using namespace std;" paragraph (they disappear entirely, along with
#    their paragraph mark, per the diff).
# ---------------------------------------------------------------------------
function Remove-ParagraphWithText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -eq ($text + "`r")) {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

[void](Remove-ParagraphWithText("#include <string>"))
[void](Remove-ParagraphWithText("using namespace std;"))

# ---------------------------------------------------------------------------
# 2. Qualify std:: members that were being used unqualified.
#    Replace every whole-word occurrence of cout/cin/string with the
#    std:: qualified form.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

[void]$find.Execute("cout", $false, $true, $false, $false, $false, $true, 1, $false, "std::cout", 2)
[void]$find.Execute("cin", $false, $true, $false, $false, $false, $true, 1, $false, "std::cin", 2)
[void]$find.Execute("string", $false, $true, $false, $false, $false, $true, 1, $false, "std::string", 2)

# ---------------------------------------------------------------------------
# 3. Switch direct field access on Linea*/Red objects to the new accessor
#    methods.
# ---------------------------------------------------------------------------
[void]$find.Execute("redMetro.cabezaLineas", $true, $false, $false, $false, $false, $true, 1, $false, "redMetro.getCabezaLineas()", 2)
[void]$find.Execute("temp->siguienteLinea", $true, $false, $false, $false, $false, $true, 1, $false, "temp->getSiguienteLinea()", 2)

# ---------------------------------------------------------------------------
# 4. Insert a new, empty paragraph right after the closing "}" of main(),
#    before the existing trailing blank paragraph.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "}`r") {
        $p.Range.InsertParagraphAfter()
        break
    }
}
